$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.664.22"
$ws.Range("E2").Value = "  -2.55%  "
$ws.Range("D3").Value = "3.408.78"
$ws.Range("E3").Value = "  -3.64%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.74"
$ws.Range("E5").Value = "  -4.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.17"
$ws.Range("E6").Value = "  -7.77%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.407.77"
$ws.Range("E8").Value = "  -3.61%  "
$ws.Range("E9").Value = "  -6.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.119"
$ws.Range("E10").Value = "  -8.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.04"
$ws.Range("E11").Value = "  -10.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.371"
$ws.Range("E12").Value = "  -9.41%  "
$ws.Range("D13").Value = "3.989.50"
$ws.Range("E13").Value = "  -3.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000176"
$ws.Range("E14").Value = "  -9.47%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "3.427.07"
$ws.Range("E15").Value = "  -3.20%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.115"
$ws.Range("E16").Value = "  -1.73%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "64.666.76"
$ws.Range("E17").Value = "  -2.40%  "
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.87"
$ws.Range("E18").Value = "  -9.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.39"
$ws.Range("E19").Value = "  -15.40%  "
$ws.Range("E20").Value = "  -7.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.38"
$ws.Range("E21").Value = "  -8.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "379.01"
$ws.Range("E22").Value = "  -10.27%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.536"
$ws.Range("E24").Value = "  -9.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.38"
$ws.Range("E25").Value = "  -7.54%  "
$ws.Range("D26").Value = "3.549.59"
$ws.Range("E26").Value = "  -3.57%  "
$ws.Range("E27").Value = "  -10.78%  "
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.10"
$ws.Range("E29").Value = "  -10.00%  "
$ws.Range("E30").Value = "  -12.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.88"
$ws.Range("E31").Value = "  -11.59%  "
$ws.Range("D32").Value = "3.429.78"
$ws.Range("E32").Value = "  -3.29%  "
$ws.Range("E34").Value = "  -9.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "22.74"
$ws.Range("E35").Value = "  -6.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "169.10"
$ws.Range("E36").Value = "  -3.68%  "
$ws.Range("E37").Value = "  -13.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.60"
$ws.Range("E38").Value = "  -13.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.43"
$ws.Range("E39").Value = "  -12.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.48"
$ws.Range("E40").Value = "  -14.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0752"
$ws.Range("E41").Value = "  -8.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.799"
$ws.Range("E42").Value = "  -7.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.83"
$ws.Range("E44").Value = "  -7.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.21"
$ws.Range("E45").Value = "  -15.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.59"
$ws.Range("E46").Value = "  -9.95%  "
$ws.Range("E47").Value = "  -1.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.20"
$ws.Range("E48").Value = "  -6.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.41"
$ws.Range("E49").Value = "  -9.27%  "
$ws.Range("D50").Value = "2.191.18"
$ws.Range("E50").Value = "  -5.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.94"
$ws.Range("E51").Value = "  -19.18%  "
